# Generate Report for Handoff
#
# This refreshes the "Latest HO Xliff Generate Date" / "Latest Handoff
# Datetime" timestamps for the 185bcd18 / 60a31766 / 80c7e31f / a34dc905 /
# c4b82cf5 / df48029c handoff rows, and marks their Priority as "ht" (hot
# fix / high priority) now that the report has been (re)generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Rows (on each sheet) that belong to the handoff batch being refreshed.
$rows = 7, 8, 11, 12, 13, 14

# "Latest HO Xliff Generate Date" on Overview (col G) -- this same
# timestamp value is also surfaced as "Latest Handoff Datetime" (col H)
# on the de-de sheet, so update both.
$hoGenerateDate = "2016-08-30 00:23:32"

# "Latest Handoff Datetime" (col H) on the zh-cn sheet.
$zhCnHandoffDate = "2016-08-30 00:23:27"

foreach ($row in $rows) {
    $wsOverview.Range("G$row").Value = $hoGenerateDate
    $wsDeDe.Range("H$row").Value     = $hoGenerateDate
    $wsZhCn.Range("H$row").Value     = $zhCnHandoffDate

    # Priority column (E) now flagged as "ht" for this handoff batch.
    $wsZhCn.Range("E$row").Value = "ht"
    $wsDeDe.Range("E$row").Value = "ht"
}
